# Update review point statuses from "Open" to "Closed" based on review.
# Rows 4, 8, 9 and 10 of the "HSI review" sheet (column E = "Point status")
# were re-reviewed and closed out.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HSI review")

$ws.Range("E4").Value = "Closed"
$ws.Range("E8").Value = "Closed"
$ws.Range("E9").Value = "Closed"
$ws.Range("E10").Value = "Closed"

# Match the author's final on-screen selection (cell E8) as recorded in
# the saved view state.
$ws.Range("E8").Select()
